$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- G13 formula change: N2 winding count 1.5x -> 2x primary ---
$ws.Range("G13").Formula = "=2*G9"

# --- Row 16: Res/m label + value (value gets a new Arial/dark-gray font) ---
$ws.Range("M16").Value = "Res/m"
$ws.Range("N16").Value = 52.96
$ws.Range("N16").Font.Color = 2237202
$ws.Range("N16").Font.Name = "Arial"

# --- Row 18: D1 ideal + F (capacitor) ---
$ws.Range("C18").Value = "D1 ideal"
$ws.Range("D18").Value = 15.65
$ws.Range("M18").Value = "F"
$ws.Range("N18").Value = 9.3

# --- Row 19: D2 ideal + M (turns) ---
$ws.Range("C19").Value = "D2 ideal"
$ws.Range("D19").Value = 31.25
$ws.Range("M19").Value = "M"
$ws.Range("N19").Value = 8

# --- Row 23: Length of winding ---
$ws.Range("M23").Value = "Length"
$ws.Range("N23").Formula = "=2*PI()*(N19+N18)/2"

# --- Row 24: D1 NON ---
$ws.Range("C24").Value = "D1 NON"
$ws.Range("D24").Value = 43.5

# --- Row 25: D2 NON + R1 ---
$ws.Range("C25").Value = "D2 NON"
$ws.Range("D25").Value = 22.05
$ws.Range("M25").Value = "R1"
$ws.Range("N25").Formula = "=N23*H9*N16/1000/2"
$ws.Range("N25").NumberFormat = "0.00E+00"

# --- Row 26: R2 ---
$ws.Range("M26").Value = "R2"
$ws.Range("N26").Formula = "=N23*N16*H13/1000/3"
$ws.Range("N26").NumberFormat = "0.00E+00"

# --- Row 27: R3 ---
$ws.Range("M27").Value = "R3"
$ws.Range("N27").Formula = "=N25"
$ws.Range("N27").NumberFormat = "0.00E+00"

# --- Row 31: Compansator section header ---
$ws.Range("D31").Value = "Compansator"

# --- Row 33-34: f_lc / f_esr ---
$ws.Range("D33").Value = "f_lc"
$ws.Range("E33").Value = 1110
$ws.Range("E33").NumberFormat = "0.00E+00"

$ws.Range("D34").Value = "f_esr"
$ws.Range("E34").Value = 1760000
$ws.Range("E34").NumberFormat = "0.00E+00"

# --- view state: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 12
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D35").Select()

# --- page setup ---
$ws.PageSetup.Orientation = 1
